$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$RangeAddr, [string]$Text)
    $helper = $ws.Range("Z1")
    $helper.NumberFormat = "@"
    $helper.Value = $Text
    $helper.Copy()
    $ws.Range($RangeAddr).PasteSpecial(-4163)
    $helper.Delete(-4159)
}

# Row 2
$ws.Range("D2").Value = "37.810.62"
$ws.Range("E2").Value = "  +1.85%  "

# Row 3
$ws.Range("D3").Value = "2.098.05"
$ws.Range("E3").Value = "  +2.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
Set-TextValue "D5" "232.88"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6
Set-TextValue "D6" "0.623"
$ws.Range("E6").Value = "  +0.52%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D8" "57.78"
$ws.Range("E8").Value = "  +1.51%  "

# Row 9
Set-TextValue "D9" "0.389"
$ws.Range("E9").Value = "  +2.05%  "

# Row 10
Set-TextValue "D10" "0.0778"
$ws.Range("E10").Value = "  +2.97%  "

# Row 11
Set-TextValue "D11" "0.104"
$ws.Range("E11").Value = "  +2.61%  "

# Row 12
$ws.Range("D12").Value = "2.395.23"
$ws.Range("E12").Value = "  +1.77%  "

# Row 13
Set-TextValue "D13" "14.54"
$ws.Range("E13").Value = "  +0.83%  "

# Row 14
Set-TextValue "D14" "21.45"
$ws.Range("E14").Value = "  +3.36%  "

# Row 15
Set-TextValue "D15" "0.778"
$ws.Range("E15").Value = "  +0.57%  "

# Row 16
Set-TextValue "D16" "5.23"
$ws.Range("E16").Value = "  +2.25%  "

# Row 17
$ws.Range("D17").Value = "2.101.56"
$ws.Range("E17").Value = "  +2.53%  "

# Row 18
$ws.Range("D18").Value = "37.792.55"
$ws.Range("E18").Value = "  +1.88%  "

# Row 19
Set-TextValue "D19" "6.18"
$ws.Range("E19").Value = "  -2.40%  "

# Row 20
Set-TextValue "D20" "70.67"
$ws.Range("E20").Value = "  +2.11%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  +2.07%  "

# Row 22
Set-TextValue "D22" "227.38"
$ws.Range("E22").Value = "  +1.21%  "

# Row 23
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("E24").Value = "  +0.62%  "

# Row 25
Set-TextValue "D25" "2.39"
$ws.Range("E25").Value = "  +0.59%  "

# Row 26
Set-TextValue "D26" "167.94"
$ws.Range("E26").Value = "  +1.46%  "

# Row 27
Set-TextValue "D27" "0.138"
$ws.Range("E27").Value = "  +9.86%  "

# Row 28
Set-TextValue "D28" "8.97"
$ws.Range("E28").Value = "  +2.52%  "

# Row 29
Set-TextValue "D29" "1.42"
$ws.Range("E29").Value = "  -1.32%  "

# Row 30
Set-TextValue "D30" "19.44"
$ws.Range("E30").Value = "  +2.53%  "

# Row 31
Set-TextValue "D31" "0.118"
$ws.Range("E31").Value = "  +1.28%  "

# Row 32
Set-TextValue "D32" "4.65"
$ws.Range("E32").Value = "  +4.74%  "

# Row 33
Set-TextValue "D33" "2.59"
$ws.Range("E33").Value = "  +3.94%  "

# Row 34
Set-TextValue "D34" "0.0624"
$ws.Range("E34").Value = "  +1.50%  "

# Row 35
Set-TextValue "D35" "4.60"
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
Set-TextValue "D36" "3.45"
$ws.Range("E36").Value = "  +5.92%  "

# Row 37
$ws.Range("E37").Value = "  +4.69%  "

# Row 38
$ws.Range("E38").Value = "  +0.07%  "

# Row 39
Set-TextValue "D39" "5.44"
$ws.Range("E39").Value = "  -4.61%  "

# Row 40
Set-TextValue "D40" "0.0983"
$ws.Range("E40").Value = "  +5.81%  "

# Row 42
Set-TextValue "D42" "96.89"
$ws.Range("E42").Value = "  +0.46%  "

# Row 43
Set-TextValue "D43" "0.0215"
$ws.Range("E43").Value = "  +2.56%  "

# Row 44
$ws.Range("D44").Value = "1.453.19"
$ws.Range("E44").Value = "  -1.69%  "

# Row 45
Set-TextValue "D45" "1.16"
$ws.Range("E45").Value = "  -1.17%  "

# Row 46
Set-TextValue "D46" "15.86"
$ws.Range("E46").Value = "  +5.36%  "

# Row 47
Set-TextValue "D47" "4.11"
$ws.Range("E47").Value = "  -6.71%  "

# Row 48
$ws.Range("E48").Value = "  +4.05%  "

# Row 49
Set-TextValue "D49" "7.28"
$ws.Range("E49").Value = "  +1.60%  "

# Row 50
Set-TextValue "D50" "3.02"
$ws.Range("E50").Value = "  +2.73%  "

# Row 51
$ws.Range("D51").Value = "2.291.66"
$ws.Range("E51").Value = "  +2.22%  "
